$d = $word.ActiveDocument

# 5.2 "clockwise" -> "anti-clockwise"
$d.Content.Find.Execute(
    "5.2 All counters are redistributed to other pits in clockwise direction.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5.2 All counters are redistributed to other pits in anti-clockwise direction.",
    2)
